# Generate Report for Handback
# - Update localization status text from "Ready for handoff" to
#   "Handed back: in sync with en-US" everywhere it appears (Overview + language sheets).
# - Refresh the "Latest Handback DateTime" for each language with the new handback
#   timestamp, and clear the stale "Error Detail" message now that the handback is
#   in sync.
# - Widen the "Status"/language status columns and shrink the "Error Detail" columns
#   to fit the new (shorter) report layout.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status cells ---
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# --- zh-cn sheet ---
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-08-16 10:49:26"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet ---
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-08-16 10:49:34"
$wsDeDe.Range("P2").Value = ""

# --- Column width adjustments ---
# NOTE: Excel's ColumnWidth setter snaps to whole-pixel increments (1/6 of a
# character unit here), so the inputs below are chosen to land on the closest
# achievable stored width to the target layout (~29.98 and ~13.75 chars).
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333334

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333334
